$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 440
$ws1.Range("F8").Value = 1992
$ws1.Range("F15").Value = 53
$ws1.Range("F19").Value = 7
$ws1.Range("F20").Value = 443
$ws1.Range("F24").Value = 7013
$ws1.Range("F25").Value = 7584
$ws1.Range("F30").Value = 77
$ws1.Range("F41").Value = 683
$ws1.Range("F44").Value = 310
$ws1.Range("F47").Value = 80
$ws1.Range("F48").Value = 127
$ws1.Range("F49").Value = 136

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 175
$ws3.Range("F4").Value = 257
$ws3.Range("F5").Value = 120

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 175
$ws4.Range("F5").Value = 120
$ws4.Range("F7").Value = 440
$ws4.Range("F8").Value = 1992
$ws4.Range("F15").Value = 53
$ws4.Range("F17").Value = 7
$ws4.Range("F18").Value = 443
$ws4.Range("F22").Value = 7013
$ws4.Range("F23").Value = 7584
$ws4.Range("F26").Value = 77
$ws4.Range("F36").Value = 683
$ws4.Range("F42").Value = 310
$ws4.Range("F45").Value = 80
$ws4.Range("F46").Value = 127
$ws4.Range("F47").Value = 136
